$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G7").Value = "CANCELLED"
$ws.Range("H11").Select()
